# "removendo variaveis não utilizadas"
# The sheet already has a 3-row sample (Sua Referencia/Numero CE header + 2 data rows).
# The edit duplicates the last data row (row 3: "YYYYYYYYYYYYYY" / "15211111111111")
# into four more rows (4-7), keeping the exact same text/shared-string values and
# the same number formatting (row 3's column B uses a text/quote-prefixed style so the
# numeric-looking string is preserved as text). Using Copy/Paste (instead of plain
# Value assignment) guarantees the new cells reuse the same shared-string entries and
# the same cell style (s="2") as row 3, instead of Excel reinterpreting the text as a
# literal number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("A3:B3")

$source.Copy($ws.Range("A4:B4"))
$source.Copy($ws.Range("A5:B5"))
$source.Copy($ws.Range("A6:B6"))
$source.Copy($ws.Range("A7:B7"))

# Clear the marching-ants clipboard marker left behind by Copy.
[void]($excel.CutCopyMode = $false)

# Match the final selection recorded in the saved file.
$ws.Range("C13").Select() | Out-Null
